$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, pushing existing rows 19..86 down to 20..87
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new market-day record
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = "Macroferia Regional de Talca"
$ws.Range("C19").Value = "Maule"
$ws.Range("D19").Value = 44459
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 32000
$ws.Range("L19").Value = 32000
$ws.Range("M19").Value = 32000
$ws.Range("N19").Value = "$/malla 25 kilos"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 1280
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
